# David Lynch Weather Tracking workbook update
#
# Adds a "Christmas note" to Sheet2 (D1) and fills in the new "From memory
# format" weather script on Sheet3 (A1:A23), then leaves Sheet3 as the
# active/selected sheet (matching the author's final view state).
#
# NOTE: cell values are written in a specific order below (not simple
# top-to-bottom row order) so that the workbook's shared-string table ends
# up with the same new-entry ordering the author's Excel session produced.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # Sheet2
$ws3 = $wb.Worksheets.Item(3)   # Sheet3

# -- Sheet2: new note in D1 --
$ws2.Range("D1").Value2 = 'on christmas, christmas message after current temps'

# -- Sheet3: new "from memory" script, written in shared-string creation order --
$ws3.Range("A1").Value2  = 'From memory format'
$ws3.Range("A2").Value2  = 'Good'
$ws3.Range("A4").Value2  = "[it's date]"
$ws3.Range("A5").Value2  = '[year]'
$ws3.Range("A7").Value2  = 'here in'
$ws3.Range("A8").Value2  = 'LA'
$ws3.Range("A11").Value2 = 'tempf'
$ws3.Range("A12").Value2 = 'degreesf'
$ws3.Range("A13").Value2 = 'tempc'
$ws3.Range("A14").Value2 = 'degreesc'
$ws3.Range("A10").Value2 = 'wind conditions'
$ws3.Range("A16").Value2 = 'later'
$ws3.Range("A17").Value2 = 'should go [up or down]'
$ws3.Range("A18").Value2 = 'forecastf'
$ws3.Range("A20").Value2 = 'forecastc'
$ws3.Range("A3").Value2  = '(Good morning or afternoon/evening)'
$ws3.Range("A6").Value2  = "and it's a [day of the week]"
$ws3.Range("A15").Value2 = 'MUSIC thinking about'
$ws3.Range("A22").Value2 = 'blue skies'
$ws3.Range("A23").Value2 = 'have a great day'

# Remaining cells reuse already-existing shared strings (order doesn't matter)
$ws3.Range("A9").Value2  = 'cloud conditions'
$ws3.Range("A19").Value2 = 'degreesf'
$ws3.Range("A21").Value2 = 'degreesc'

# -- View state: Sheet3 becomes the selected/active tab --
$ws2.Activate()
$ws2.Range("H5").Select()
$ws3.Activate()
$ws3.Range("F16").Select()
